$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44301, 1, 11, 118.4834123222749),
    @(44302, 3, 14, 150.7970702283498),
    @(44303, 1, 9, 96.9409737182249)
)

$startRow = 227
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

# Copy the date-column formatting (style used for column A) from the last existing
# row down onto the newly added rows, matching the sheet's existing style pattern.
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A${startRow}:A$($startRow + $newRows.Count - 1)").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
